# Update "想去人数" (F column) counts on sheets "展览" and "全部类型".
# Both sheets contain identical data tables, so the same row/value
# changes are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 866
    5  = 65
    6  = 12707
    7  = 62
    8  = 108
    10 = 463
    13 = 13676
    14 = 14053
    16 = 167
    19 = 18
    23 = 1061
    25 = 56
    26 = 930
    27 = 5172
    29 = 261
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
